$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    # Force the cell to stay text-typed (matches source data, which stores
    # numeric-looking prices/percentages as literal strings) instead of Excel
    # auto-coercing to a Number/Percentage, then strip the format stamp back off
    # so the style index is unchanged from the original (unstyled) cell.
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.ClearFormats()
}

Set-TextValue "D2" '289.92'
Set-TextValue "E2" '-3.99%'
Set-TextValue "D3" '30.71'
Set-TextValue "E3" '-4.76%'
Set-TextValue "D4" '4.883'
Set-TextValue "E4" '-1.87%'
Set-TextValue "D5" '0.07203'
Set-TextValue "E5" '-8.87%'
$ws.Range("B6").Value = 'FTXToken'
$ws.Range("C6").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-TextValue "D6" '1.779'
Set-TextValue "E6" '-16.25%'
$ws.Range("B7").Value = 'KuCoinToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
Set-TextValue "D7" '7.675'
Set-TextValue "E7" '-2.49%'
Set-TextValue "D8" '3.735'
Set-TextValue "E8" '-1.73%'
Set-TextValue "D9" '0.8956'
Set-TextValue "E9" '-3.18%'
Set-TextValue "D10" '0.1652'
Set-TextValue "E10" '-5.66%'
Set-TextValue "D11" '0.07470'
Set-TextValue "E11" '-6.26%'
Set-TextValue "D12" '0.08054'
Set-TextValue "E12" '-7.62%'
Set-TextValue "D13" '0.02990'
Set-TextValue "E13" '-4.65%'
Set-TextValue "D14" '0.09994'
Set-TextValue "E14" '-0.24%'
Set-TextValue "D15" '0.001496'
Set-TextValue "E15" '-0.83%'
Set-TextValue "D16" '0.005681'
Set-TextValue "E16" '-1.22%'
Set-TextValue "D19" '2.100'
Set-TextValue "E19" '-7.74%'
Set-TextValue "D22" '4.396'
Set-TextValue "E22" '1.39%'
Set-TextValue "D24" '0.04479'
Set-TextValue "E24" '-2.65%'
Set-TextValue "D25" '0.001211'
Set-TextValue "E25" '-2.16%'
Set-TextValue "D26" '0.004022'
Set-TextValue "E26" '-9.90%'
Set-TextValue "D27" '0.0001251'
Set-TextValue "E27" '-0.06%'
Set-TextValue "D39" '0.01646'
Set-TextValue "E39" '-4.72%'
Set-TextValue "D40" '0.04339'
Set-TextValue "E40" '-9.47%'
Set-TextValue "D41" '0.007419'
Set-TextValue "E41" '-1.04%'
Set-TextValue "D42" '0.1311'
Set-TextValue "E42" '-3.59%'
Set-TextValue "D43" '0.002008'
Set-TextValue "E43" '-14.32%'
Set-TextValue "D44" '0.01018'
Set-TextValue "E44" '-0.65%'
Set-TextValue "D45" '0.00005809'
Set-TextValue "E45" '-3.22%'
Set-TextValue "E46" '0.09%'
Set-TextValue "D47" '2.193'
Set-TextValue "E47" '167.35%'
Set-TextValue "D48" '0.003003'
Set-TextValue "E48" '-11.43%'
Set-TextValue "D49" '0.00002102'
Set-TextValue "E49" '0.09%'
Set-TextValue "D50" '0.0002002'
Set-TextValue "E50" '0.09%'
